$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Periodo Mora" values in column E (rows 16-22) into ascending
# order (they are stored as text, e.g. "1607", "1608", ...).
$periods = @("1607", "1608", "1609", "1610", "1611", "1612", "1701")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}
